$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Sheet1" to "Excel"
$ws.Name = "Excel"

# Move selection from C7 to A2
$ws.Range("A2").Select()
